$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 relabeled to GNN-MT-O, new value
$ws.Range("A2").Value = "8_train (GNN-MT-O) val delta-auprc"
$ws.Range("B2").Value = 3.488235294117647

# Row 3 stays RF, new value
$ws.Range("A3").Value = "8_train (RF) val delta-auprc"
$ws.Range("B3").Value = 3.267676767676768

# Row 4 stays GNN-MT, new value
$ws.Range("A4").Value = "8_train (GNN-MT) val delta-auprc"
$ws.Range("B4").Value = 3.085294117647059

# New row 5: PN-O, copy formatting from row 4's label cell
$ws.Range("A5").Value = "8_train (PN-O) val delta-auprc"
$ws.Range("B5").Value = 2.441919191919192
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# New row 6: PN, copy formatting from row 4's label cell
$ws.Range("A6").Value = "8_train (PN) val delta-auprc"
$ws.Range("B6").Value = 2.373737373737374
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
